$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(5, 3).Value = 4.42
$ws.Cells.Item(5, 4).Value = "Farmace"
$ws.Cells.Item(5, 5).Value = "Farmace"
$ws.Cells.Item(6, 3).Value = 7.15
$ws.Cells.Item(6, 4).Value = "Hipolabor"
$ws.Cells.Item(6, 5).Value = "Hipolabor"
$ws.Cells.Item(7, 3).Value = 2.99
$ws.Cells.Item(7, 4).Value = "Hypofarma"
$ws.Cells.Item(7, 5).Value = "Hypofarma"
$ws.Cells.Item(8, 3).Value = 0.65
$ws.Cells.Item(8, 4).Value = "Isofarma"
$ws.Cells.Item(8, 5).Value = "Isofarma"
$ws.Cells.Item(10, 3).Value = 12.87
$ws.Cells.Item(10, 4).Value = "Hipolabor"
$ws.Cells.Item(10, 5).Value = "Hipolabor"
$ws.Cells.Item(11, 3).Value = 76.7
$ws.Cells.Item(11, 4).Value = "EMS"
$ws.Cells.Item(11, 5).Value = "EMS"
$ws.Cells.Item(21, 3).Value = 12.35
$ws.Cells.Item(21, 4).Value = "Teuto"
$ws.Cells.Item(21, 5).Value = "Teuto"
$ws.Cells.Item(27, 3).Value = 3.38
$ws.Cells.Item(27, 4).Value = "Hypofarma"
$ws.Cells.Item(27, 5).Value = "Hypofarma"
$ws.Cells.Item(35, 3).Value = 8.97
$ws.Cells.Item(35, 4).Value = "Blau"
$ws.Cells.Item(35, 5).Value = "Blau"
$ws.Cells.Item(36, 3).Value = 17.55
$ws.Cells.Item(36, 4).Value = "Teuto"
$ws.Cells.Item(36, 5).Value = "Teuto"
$ws.Cells.Item(37, 3).Value = 6.37
$ws.Cells.Item(37, 4).Value = "Blau"
$ws.Cells.Item(37, 5).Value = "Blau"
$ws.Cells.Item(41, 3).Value = 3.25
$ws.Cells.Item(41, 4).Value = "Teuto"
$ws.Cells.Item(41, 5).Value = "Teuto"
$ws.Cells.Item(43, 3).Value = 44.2
$ws.Cells.Item(43, 4).Value = "Isofarma"
$ws.Cells.Item(43, 5).Value = "Isofarma"
$ws.Cells.Item(63, 3).Value = 3.9
$ws.Cells.Item(63, 4).Value = "Sanval"
$ws.Cells.Item(63, 5).Value = "Sanval"
$ws.Cells.Item(64, 3).Value = 10.4
$ws.Cells.Item(64, 4).Value = "Cristália"
$ws.Cells.Item(64, 5).Value = "Cristália"
$ws.Cells.Item(66, 3).Value = 3.51
$ws.Cells.Item(66, 4).Value = "Teuto"
$ws.Cells.Item(66, 5).Value = "Teuto"
$ws.Cells.Item(71, 3).Value = 1.95
$ws.Cells.Item(71, 4).Value = "Teuto"
$ws.Cells.Item(71, 5).Value = "Teuto"
$ws.Cells.Item(73, 3).Value = 24.7
$ws.Cells.Item(73, 4).Value = "Takeda"
$ws.Cells.Item(73, 5).Value = "Takeda"
$ws.Cells.Item(75, 3).Value = 2.08
$ws.Cells.Item(75, 4).Value = "Teuto"
$ws.Cells.Item(75, 5).Value = "Teuto"
$ws.Cells.Item(81, 3).Value = 89.7
$ws.Cells.Item(81, 4).Value = "Cristália"
$ws.Cells.Item(81, 5).Value = "Cristália"
$ws.Cells.Item(82, 3).Value = 29.9
$ws.Cells.Item(82, 4).Value = "Mylan"
$ws.Cells.Item(82, 5).Value = "Mylan"
$ws.Cells.Item(92, 3).Value = 7.15
$ws.Cells.Item(92, 4).Value = "Hypofarma"
$ws.Cells.Item(92, 5).Value = "Hypofarma"
$ws.Cells.Item(94, 3).Value = 2.47
$ws.Cells.Item(94, 4).Value = "Teuto"
$ws.Cells.Item(94, 5).Value = "Teuto"
$ws.Cells.Item(100, 3).Value = 1.24
$ws.Cells.Item(100, 4).Value = "Isofarma"
$ws.Cells.Item(100, 5).Value = "Isofarma"
$ws.Cells.Item(103, 3).Value = 11.57
$ws.Cells.Item(103, 4).Value = "Cristália"
$ws.Cells.Item(103, 5).Value = "Cristália"
$ws.Cells.Item(107, 3).Value = 6.37
$ws.Cells.Item(107, 4).Value = "Blau"
$ws.Cells.Item(107, 5).Value = "Blau"
$ws.Cells.Item(108, 3).Value = 9.09
$ws.Cells.Item(108, 4).Value = "Teuto"
$ws.Cells.Item(108, 5).Value = "Teuto"
$ws.Cells.Item(123, 3).Value = 8.97
$ws.Cells.Item(123, 4).Value = "Halex Istar"
$ws.Cells.Item(123, 5).Value = "Halex Istar"
$ws.Cells.Item(134, 3).Value = 12.09
$ws.Cells.Item(134, 4).Value = "Blau"
$ws.Cells.Item(134, 5).Value = "Blau"
$ws.Cells.Item(136, 3).Value = 3.77
$ws.Cells.Item(136, 4).Value = "Hipolabor"
$ws.Cells.Item(136, 5).Value = "Hipolabor"
$ws.Cells.Item(142, 3).Value = 36.4
$ws.Cells.Item(142, 4).Value = "Mylan"
$ws.Cells.Item(142, 5).Value = "Mylan"
$ws.Cells.Item(143, 3).Value = 3.89
$ws.Cells.Item(143, 4).Value = "Hypofarma"
$ws.Cells.Item(143, 5).Value = "Hypofarma"
$ws.Cells.Item(144, 3).Value = 20.8
$ws.Cells.Item(144, 4).Value = "Cristália"
$ws.Cells.Item(144, 5).Value = "Cristália"
$ws.Cells.Item(145, 3).Value = 4
$ws.Cells.Item(145, 4).Value = "Belfar"
$ws.Cells.Item(145, 5).Value = "Belfar"
$ws.Cells.Item(146, 3).Value = 14.82
$ws.Cells.Item(146, 4).Value = "Nativita"
$ws.Cells.Item(146, 5).Value = "Nativita"
$ws.Cells.Item(155, 3).Value = 12.35
$ws.Cells.Item(155, 4).Value = "JP"
$ws.Cells.Item(155, 5).Value = "JP"
$ws.Cells.Item(156, 3).Value = 8.32
$ws.Cells.Item(156, 4).Value = "JP"
$ws.Cells.Item(156, 5).Value = "JP"
$ws.Cells.Item(157, 3).Value = 9.36
$ws.Cells.Item(157, 4).Value = "JP"
$ws.Cells.Item(157, 5).Value = "JP"
$ws.Cells.Item(158, 3).Value = 9.88
$ws.Cells.Item(158, 4).Value = "Halex Istar"
$ws.Cells.Item(158, 5).Value = "Halex Istar"
$ws.Cells.Item(159, 3).Value = 14.17
$ws.Cells.Item(159, 4).Value = "JP"
$ws.Cells.Item(159, 5).Value = "JP"
$ws.Cells.Item(160, 3).Value = 12.87
$ws.Cells.Item(160, 4).Value = "Halex Istar"
$ws.Cells.Item(160, 5).Value = "Halex Istar"
$ws.Cells.Item(163, 3).Value = 12.87
$ws.Cells.Item(163, 4).Value = "JP"
$ws.Cells.Item(163, 5).Value = "JP"
$ws.Cells.Item(164, 3).Value = 3.77
$ws.Cells.Item(164, 4).Value = "Farmace"
$ws.Cells.Item(164, 5).Value = "Farmace"
$ws.Cells.Item(170, 3).Value = 21.48
$ws.Cells.Item(170, 4).Value = "Teuto"
$ws.Cells.Item(170, 5).Value = "Teuto"
$ws.Cells.Item(172, 3).Value = 11.7
$ws.Cells.Item(172, 4).Value = "Cristália"
$ws.Cells.Item(172, 5).Value = "Cristália"
$ws.Cells.Item(174, 3).Value = 3.25
$ws.Cells.Item(174, 4).Value = "Teuto"
$ws.Cells.Item(174, 5).Value = "Teuto"
$ws.Cells.Item(183, 3).Value = 2.08
$ws.Cells.Item(183, 4).Value = "Teuto"
$ws.Cells.Item(183, 5).Value = "Teuto"
$ws.Cells.Item(185, 3).Value = 23.4
$ws.Cells.Item(185, 4).Value = "Blau"
$ws.Cells.Item(185, 5).Value = "Blau"
$ws.Cells.Item(191, 3).Value = 14.3
$ws.Cells.Item(191, 4).Value = "Cristália"
$ws.Cells.Item(191, 5).Value = "Cristália"
$ws.Cells.Item(192, 3).Value = 7.67
$ws.Cells.Item(192, 4).Value = "União Química"
$ws.Cells.Item(192, 5).Value = "União Química"
$ws.Cells.Item(196, 3).Value = 11.7
$ws.Cells.Item(196, 4).Value = "Eurofarma"
$ws.Cells.Item(196, 5).Value = "Eurofarma"
